$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "[채용] IT 기획 담당자 찾습니다 (대학원 Data Analytics 과정 장학생 추천)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/pabii-hr-20210304/#utm_source=rss&utm_medium=rss&utm_campaign=pabii-hr-20210304"

$ws.Range("D32").Value = "Autoencoder 를 이용한 차원 축소 (latent representation)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/301"

$ws.Range("D39").Value = "Apparent Age and Gender Prediction in Keras"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Apparent-Age-and-Gender-Prediction-in-Keras-1"

$ws.Range("D45").Value = "Time-series forecasting at UBER"
$ws.Range("E45").Value = "https://dive-into-ds.tistory.com/83"

$ws.Range("D51").Value = "[세이버메트릭스] 평균 연봉이 더 높은 구단이 더 잘했을까? (2019년, 2020년 분석)"
$ws.Range("E51").Value = "https://bskyvision.com/1127"
